# Atualização de bases das ligas, do dia: 14-06-2024 às 20:31
#
# The underlying source rows did not change; several match rows (which all
# share the same Div/Date) had their data (id/B through AD, i.e. everything
# except columns A, C and D) reassigned among each other. This script
# reproduces that permutation by snapshotting the affected rows' values
# first and then writing them to their new positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-RowGroup {
    param(
        [object]$Worksheet,
        [int[]]$RowOrder
    )

    # Snapshot column B (id/match key) and columns E:AD (everything else
    # that moves) for every row in the group before writing anything,
    # so rows that are sources for one destination are not clobbered
    # before they have been read.
    $countRows = $RowOrder.Count
    $snapB = @{}
    $snapEAD = @{}
    foreach ($r in $RowOrder) {
        $snapB[$r] = $Worksheet.Range("B$r").Value()
        $snapEAD[$r] = $Worksheet.Range("E${r}:AD$r").Value()
    }

    # Row at index i receives the snapshot of the row at index i+1
    # (cyclically). This matches the permutation found in the diff.
    for ($i = 0; $i -lt $countRows; $i++) {
        $dst = $RowOrder[$i]
        $srcRow = $RowOrder[($i + 1) % $countRows]
        $Worksheet.Range("B$dst").Value = $snapB[$srcRow]
        $Worksheet.Range("E${dst}:AD$dst").Value = $snapEAD[$srcRow]
    }
}

Copy-RowGroup $ws @(42, 43, 44, 45)
Copy-RowGroup $ws @(49, 50)
Copy-RowGroup $ws @(67, 69, 68)
Copy-RowGroup $ws @(97, 98)
Copy-RowGroup $ws @(112, 114, 115)
Copy-RowGroup $ws @(120, 122, 121)
Copy-RowGroup $ws @(124, 125)
